$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126; this shifts existing rows 126:221 down to 127:222
$ws.Rows(126).Insert()

# Populate the newly inserted row 126 with the new weekly data point
$ws.Range("A126").Value = 8
$ws.Range("B126").Value = "Terminal La Palmera de La Serena"
$ws.Range("C126").Value = "Coquimbo"
$ws.Range("D126").Value = 44447
$ws.Range("E126").Value = 4
$ws.Range("F126").Value = 100114001
$ws.Range("G126").Value = "Papa"
$ws.Range("H126").Value = "Cardinal"
$ws.Range("I126").Value = "1a (cosecha)"
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 11000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = 11500
$ws.Range("N126").Value = "`$/saco 25 kilos"
$ws.Range("O126").Value = "Provincia del Elquí"
$ws.Range("P126").Value = 460
$ws.Range("Q126").Value = 25
$ws.Range("R126").Value = "Hortaliza"
